$wb = $excel.ActiveWorkbook

# Sheet 1: VENTAS POR GRUPO
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("I14").Value = 463.5
$ws1.Range("L14").Value = 1690.58
$ws1.Range("M14").Value = 12167.91
$ws1.Range("I55").Value = "6 de 53"

# Sheet 2: VENTA MENSUAL
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F14").Value = 15019.68
$ws2.Range("F55").Value = 100865.71

# Sheet 3: CUMPLIMIENTO MENSUAL
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D8").Value = 1891.59
$ws3.Range("E8").Value = -891.5899999999999
$ws3.Range("F8").Value = 1.89159

$ws3.Range("D15").Value = 7758.83
$ws3.Range("E15").Value = 5741.17
$ws3.Range("F15").Value = 0.5747281481481481

$ws3.Range("D16").Value = 53629.42
$ws3.Range("E16").Value = 2430.279999999999
$ws3.Range("F16").Value = 0.9566483588032044

$ws3.Range("D19").Value = 100865.71
$ws3.Range("E19").Value = 16573.98064517915
$ws3.Range("F19").Value = 0.8588724088583121
